# Insert a new weekly price-record row before row 123 (Mango, Macroferia
# Regional de Talca), shifting the existing rows 123-153 down to 124-154.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(123).Insert()

$ws.Range("A123").Value = 5
$ws.Range("B123").Value = "Macroferia Regional de Talca"
$ws.Range("C123").Value = "Maule"
$ws.Range("D123").Value = 44855
$ws.Range("E123").Value = 7
$ws.Range("F123").Value = "Fruta"
$ws.Range("G123").Value = 100108
$ws.Range("H123").Value = "Tropicales y subtropicales"
$ws.Range("I123").Value = 100108002
$ws.Range("J123").Value = "Mango"
$ws.Range("K123").Value = "Sin especificar"
$ws.Range("L123").Value = "Primera"
$ws.Range("M123").Value = 240
$ws.Range("N123").Value = 8000
$ws.Range("O123").Value = 8000
$ws.Range("P123").Value = 8000
$ws.Range("Q123").Value = "`$/bandeja 4 kilos"
$ws.Range("R123").Value = "Brasil"
$ws.Range("S123").Value = 2000
$ws.Range("T123").Value = 4
